$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected Neo4j query text for the "dbExcel" query cell (B2):
# "WHERE c.ethnicity IN ['UNKNOWN']" -> "WHERE c.ethnicity = "UNKNOWN" "
$query1 = @'
MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)
    WHERE c.ethnicity = "UNKNOWN" 
WITH DISTINCT c, a, ct
RETURN 
    COALESCE(c.case_id, '') AS `Case ID`,
    COALESCE(ct.clinical_trial_designation, '') AS `Trial Code`,
    COALESCE(a.arm_id, '') AS `Arm`,
    COALESCE(a.arm_drug, '') AS `Arm Treatment`,
    COALESCE(c.disease, '') AS `Diagnosis`,
    COALESCE(c.gender, '') AS `Gender`,
    COALESCE(c.race, '') AS `Race`,
    COALESCE(c.ethnicity, '') AS `Ethnicity`
'@

# Corrected Neo4j query text for the "WebExcel" stat query cell (C2):
# "WHERE WHERE c.ethnicity IN ['UNKNOWN']" -> "WHERE c.ethnicity = "UNKNOWN" "
$query2 = @'
MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE c.ethnicity = "UNKNOWN" 
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
    COUNT(DISTINCT f) AS number_of_files,
    COUNT(DISTINCT c.case_id) AS number_of_cases,
    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials
'@

$ws.Range("B2").Value = $query1
$ws.Range("C2").Value = $query2

# Update the active selection to match the saved view state (B3)
$ws.Range("B3").Select()
